# ---------------------------------------------------------------------------
# feat(publipostage): Add "status_label" as string version of "status"
#
# The sheet gains a new column B ("status_label") holding a plain-text French
# word (rouge / orange) describing the colored-circle emoji already stored in
# column A (statut). The former columns B:I (NCTId .. results) are pushed one
# column to the right, becoming C:J.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank column before column B. Excel shifts the existing
#    NCTId..results data (columns B:I) one position to the right (C:J) and
#    grows the sheet dimension from I26 to J26 automatically.
$ws.Columns.Item(2).Insert()

# 2) Header for the freshly inserted column.
$ws.Range("B1").Value = 'status_label'

# 3) status_label value for every data row, derived from the emoji in column A
#    (🟥 -> rouge, 🟧 -> orange).
$ws.Range("B2").Value = 'rouge'
$ws.Range("B3").Value = 'rouge'
$ws.Range("B4").Value = 'rouge'
$ws.Range("B5").Value = 'rouge'
$ws.Range("B6").Value = 'rouge'
$ws.Range("B7").Value = 'orange'
$ws.Range("B8").Value = 'rouge'
$ws.Range("B9").Value = 'rouge'
$ws.Range("B10").Value = 'rouge'
$ws.Range("B11").Value = 'rouge'
$ws.Range("B12").Value = 'rouge'
$ws.Range("B13").Value = 'rouge'
$ws.Range("B14").Value = 'orange'
$ws.Range("B15").Value = 'rouge'
$ws.Range("B16").Value = 'rouge'
$ws.Range("B17").Value = 'orange'
$ws.Range("B18").Value = 'rouge'
$ws.Range("B19").Value = 'orange'
$ws.Range("B20").Value = 'orange'
$ws.Range("B21").Value = 'rouge'
$ws.Range("B22").Value = 'rouge'
$ws.Range("B23").Value = 'rouge'
$ws.Range("B24").Value = 'rouge'
$ws.Range("B25").Value = 'rouge'
$ws.Range("B26").Value = 'rouge'

# 4) The source export also re-ordered a handful of rows that share the same
#    completion_year. Re-write those rows fully (columns A and C:J) so the
#    final row order/content matches the published sheet exactly.

# Row 7: NCTId = NCT00862329
$ws.Range("A7").Value = '🟧'
$ws.Range("C7").Value = 'NCT00862329'
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = '2010'
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = $false
$ws.Range("I7").Value = $true
$ws.Range("J7").Value = $true

# Row 8: NCTId = NCT00690781
$ws.Range("A8").Value = '🟥'
$ws.Range("C8").Value = 'NCT00690781'
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = '2010'
$ws.Range("F8").Value = 'Effect of Milk Proteins and Protein Feeding Pattern on Body Composition and Protein Metabolism in Energy Restricted Obese Subjects'
$ws.Range("G8").Value = 'SURPROL-CF-H'
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = $false
$ws.Range("J8").Value = $false

# Row 9: NCTId = NCT01209572
$ws.Range("A9").Value = '🟥'
$ws.Range("C9").Value = 'NCT01209572'
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = '2010'
$ws.Range("F9").Value = 'Modelling of 24h Energy Expenditure From Heart Rate, Actimetry and Other Parameters Recorded Under Free-living Conditions'
$ws.Range("G9").Value = 'Modelheart'
$ws.Range("H9").Value = $false
$ws.Range("I9").Value = $false
$ws.Range("J9").Value = $false

# Row 13: NCTId = NCT02473302
$ws.Range("A13").Value = '🟥'
$ws.Range("C13").Value = 'NCT02473302'
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = '2014'
$ws.Range("F13").Value = 'Preventive Strategies in Colorectal Carcinogenesis Production and Meat Processing'
$ws.Range("G13").ClearContents()
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("J13").Value = $false

# Row 14: NCTId = NCT02354794
$ws.Range("A14").Value = '🟧'
$ws.Range("C14").Value = 'NCT02354794'
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = '2014'
$ws.Range("F14").Value = 'Effect of Oral Supplementation With One Form of L-arginine on Vascular Endothelial Function in Healthy Subjects Featuring Risk Factors Related to the Metabolic Syndrome.'
$ws.Range("G14").ClearContents()
$ws.Range("H14").Value = $false
$ws.Range("I14").Value = $true
$ws.Range("J14").Value = $true

# Row 15: NCTId = NCT02157805
$ws.Range("A15").Value = '🟥'
$ws.Range("C15").Value = 'NCT02157805'
$ws.Range("D15").ClearContents()
$ws.Range("E15").Value = '2014'
$ws.Range("F15").Value = 'Effect of Technological Processes on Nutritional Quality of Meat Proteins'
$ws.Range("G15").ClearContents()
$ws.Range("H15").Value = $false
$ws.Range("I15").Value = $false
$ws.Range("J15").Value = $true

# Row 21: NCTId = NCT06624033
$ws.Range("A21").Value = '🟥'
$ws.Range("C21").Value = 'NCT06624033'
$ws.Range("D21").ClearContents()
$ws.Range("E21").Value = '2023'
$ws.Range("F21").Value = 'Single-blind, Randomized, Cross-over Comparative Bioavailability Study About the Kinetics of Plasma Amino Acid Concentrations Subsequent to the Consumption of Innovative Legume-based Products.'
$ws.Range("G21").Value = 'LEG''UP'
$ws.Range("H21").Value = $false
$ws.Range("I21").Value = $false
$ws.Range("J21").Value = $false

# Row 22: NCTId = NCT05047757
$ws.Range("A22").Value = '🟥'
$ws.Range("C22").Value = 'NCT05047757'
$ws.Range("D22").ClearContents()
$ws.Range("E22").Value = '2023'
$ws.Range("F22").Value = 'Fava Bean Protein and Amino Acid Bioavailability in Healthy Volunteers'
$ws.Range("G22").Value = 'Leg4Life'
$ws.Range("H22").Value = $false
$ws.Range("I22").Value = $false
$ws.Range("J22").Value = $false

